# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newNote = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.54 = 50765.94 pesos`n✅ 50765.94 pesos = 12.55 = 980.15 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newNote

# --- tasas: update the N10/O10/N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 79.71899999999999
$wsTasas.Range("O10").Value = 4047.01

$wsTasas.Range("N12").Value = 4045
$wsTasas.Range("O12").Value = 78.098
